$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels (row 2 across all 3 repeated blocks: B-I, J-Q, R-Y) ---
# Column mapping before -> after (within each 8-col block):
#   col1 (e.g. D): "IMF - Sales"          -> "IMF (20%) - Sales"
#   col2 (e.g. E): "IMF - Sales + Emp"     -> "IMF (20%) - Sales + Emp"
#   col3 (e.g. F): "OECD (20%) - Sales"    -> "IMF - Sales"
#   col4 (e.g. G): "OECD (20%) - Sales + Emp" -> "IMF - Sales + Emp"
# col5-8 (OECD - Sales / OECD - Sales + Emp, repeated headers) stay the same.

$ws.Range("D2").Value = "IMF (20%) - Sales"
$ws.Range("E2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("F2").Value = "IMF - Sales"
$ws.Range("G2").Value = "IMF - Sales + Emp"

$ws.Range("L2").Value = "IMF (20%) - Sales"
$ws.Range("M2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("N2").Value = "IMF - Sales"
$ws.Range("O2").Value = "IMF - Sales + Emp"

$ws.Range("T2").Value = "IMF (20%) - Sales"
$ws.Range("U2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("V2").Value = "IMF - Sales"
$ws.Range("W2").Value = "IMF - Sales + Emp"

# --- Update data cells (rows 4-13) to match the new column layout ---
# Row 4
$ws.Range("D4").Value = 0.7733813234692754
$ws.Range("E4").Value = 0.8246446903403367
$ws.Range("F4").Value = 3.866906617346374
$ws.Range("G4").Value = 4.123223451701683
$ws.Range("N4").Value = 0.275684440409188
$ws.Range("O4").Value = 0.2753197867554495
$ws.Range("V4").Value = 62346466255
$ws.Range("W4").Value = 62520144000

# Row 5
$ws.Range("D5").Value = 0.7383447073709195
$ws.Range("E5").Value = 0.8065494146227172
$ws.Range("F5").Value = 3.691723536854596
$ws.Range("G5").Value = 4.032747073113593
$ws.Range("N5").Value = 0.2187022724763601
$ws.Range("O5").Value = 0.2201697804353018
$ws.Range("V5").Value = 976531986457
$ws.Range("W5").Value = 988562844368

# Row 6
$ws.Range("D6").Value = 0.1800287566322612
$ws.Range("E6").Value = 0.5127070547825061
$ws.Range("F6").Value = 0.9001437831613054
$ws.Range("G6").Value = 2.56353527391253
$ws.Range("N6").Value = 0.3972232104675282
$ws.Range("O6").Value = 0.3787799051437756
$ws.Range("V6").Value = 25762595315
$ws.Range("W6").Value = 37793453226

# Row 7
$ws.Range("D7").Value = 0.2532044456230912
$ws.Range("E7").Value = 0.7184075403248512
$ws.Range("F7").Value = 1.266022228115456
$ws.Range("G7").Value = 3.592037701624256
$ws.Range("N7").Value = 0.4812503551453929
$ws.Range("O7").Value = 0.4369507789038161
$ws.Range("V7").Value = 35956611724
$ws.Range("W7").Value = 49982592968

# Row 8
$ws.Range("D8").Value = 1.260193722463467
$ws.Range("E8").Value = 0.9921647696433968
$ws.Range("F8").Value = 6.300968612317333
$ws.Range("G8").Value = 4.960823848216997
$ws.Range("N8").Value = 0.2083789632708606
$ws.Range("O8").Value = 0.2083789632708606
$ws.Range("V8").Value = 933471841988
$ws.Range("W8").Value = 933471841988

# Row 9
$ws.Range("D9").Value = 0.9736807634004252
$ws.Range("E9").Value = 1.808176613332981
$ws.Range("F9").Value = 4.868403817002125
$ws.Range("G9").Value = 9.040883066664906
$ws.Range("N9").Value = 0.7565071660679659
$ws.Range("O9").Value = 0.4505402927164314
$ws.Range("V9").Value = 11068720584
$ws.Range("W9").Value = 31122850519

# Row 10
$ws.Range("D10").Value = 0.9868244542315268
$ws.Range("E10").Value = 1.073704801281935
$ws.Range("F10").Value = 4.934122271157631
$ws.Range("G10").Value = 5.368524006409673
$ws.Range("N10").Value = 0.2733403708402601
$ws.Range("O10").Value = 0.2733403708402601
$ws.Range("V10").Value = 65619795685
$ws.Range("W10").Value = 65619795685

# Row 11
$ws.Range("D11").Value = 1.056821244091495
$ws.Range("E11").Value = 0.9597021300649335
$ws.Range("F11").Value = 5.284106220457476
$ws.Range("G11").Value = 4.798510650324678
$ws.Range("N11").Value = 0.2159819121733793
$ws.Range("O11").Value = 0.2159690554647536
$ws.Range("V11").Value = 966574409108
$ws.Range("W11").Value = 966748086853

# Row 12
$ws.Range("D12").Value = 0.7943293633077065
$ws.Range("E12").Value = 0.9810512255814401
$ws.Range("F12").Value = 3.971646816538531
$ws.Range("G12").Value = 4.905256127907198
$ws.Range("N12").Value = 0.3259974015367673
$ws.Range("O12").Value = 0.3256785167725302
$ws.Range("V12").Value = 98722362805
$ws.Range("W12").Value = 98896040550

# Row 13 (only D,E,F,G change; N,O,V,W already equal the values they'd be reassigned to)
$ws.Range("D13").Value = 1.486349685879603
$ws.Range("E13").Value = 0.9247074269381814
$ws.Range("F13").Value = 7.431748429398017
$ws.Range("G13").Value = 4.623537134690935

